# Daily attendance processing - 2026-01-04 09:56:54
# Normalizes the "Recorded By" (column G) values: for any entry that is a
# comma-separated list of recorders, the first and last entries are swapped
# (moving the leading "System" marker to the end of the list), except for
# entries that include "admin@admin.com", which are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) { continue }
    if ($val -like "*admin@admin.com*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -gt 1) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $newVal = [string]::Join(", ", $parts)
        $cell.Value = $newVal
    }
}
